$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.058.24'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.883.87'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9983'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9980'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4968'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '44.40'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.84%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06618'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.880.63'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.82'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.07193'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6641'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '85.74'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.845'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.028.43'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007763'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +3.65%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.9985'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.75'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.120.51'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9976'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.760'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.602'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.152'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '150.53'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +4.39%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '136.26'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.26%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '16.78'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.379'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.165'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.08679'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.945'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04991'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -1.33%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7025'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.658'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.67%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.697'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.196'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9339'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -2.85%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +1.17%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.956'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9992'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4191'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '101.25'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.493'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.05717'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.69%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.260'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.343'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.77%  '
